$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to remain plain text so values like
# "43.858.24" or "7.90" are not re-interpreted as numbers/dates by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "43.858.24"
$ws.Range("E2").Value = "  -0.89%  "

# Row 3
$ws.Range("D3").Value = "2.308.43"
$ws.Range("E3").Value = "  +2.13%  "

# Row 4
$ws.Range("E4").Value = "  +0.21%  "

# Row 5
$ws.Range("D5").Value = "95.93"
$ws.Range("E5").Value = "  -2.94%  "

# Row 6
$ws.Range("D6").Value = "268.63"
$ws.Range("E6").Value = "  -1.58%  "

# Row 7
$ws.Range("D7").Value = "0.627"
$ws.Range("E7").Value = "  -0.91%  "

# Row 8
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
$ws.Range("D9").Value = "0.615"
$ws.Range("E9").Value = "  -5.45%  "

# Row 10
$ws.Range("D10").Value = "44.47"
$ws.Range("E10").Value = "  -7.33%  "

# Row 11
$ws.Range("D11").Value = "0.0945"
$ws.Range("E11").Value = "  -1.30%  "

# Row 12
$ws.Range("D12").Value = "7.90"
$ws.Range("E12").Value = "  -6.99%  "

# Row 13
$ws.Range("D13").Value = "0.105"
$ws.Range("E13").Value = "  -0.15%  "

# Row 14
$ws.Range("D14").Value = "2.654.19"
$ws.Range("E14").Value = "  +2.34%  "

# Row 15
$ws.Range("D15").Value = "15.33"
$ws.Range("E15").Value = "  -1.25%  "

# Row 16
$ws.Range("D16").Value = "0.859"
$ws.Range("E16").Value = "  +3.02%  "

# Row 17
$ws.Range("D17").Value = "2.314.29"
$ws.Range("E17").Value = "  +2.36%  "

# Row 18
$ws.Range("D18").Value = "43.825.59"
$ws.Range("E18").Value = "  -0.83%  "

# Row 19
$ws.Range("D19").Value = "0.0000109"
$ws.Range("E19").Value = "  +1.80%  "

# Row 20
$ws.Range("D20").Value = "6.28"
$ws.Range("E20").Value = "  -0.12%  "

# Row 21
$ws.Range("D21").Value = "73.28"
$ws.Range("E21").Value = "  +2.47%  "

# Row 22
$ws.Range("D22").Value = "2.35"
$ws.Range("E22").Value = "  +0.67%  "

# Row 23
$ws.Range("D23").Value = "238.14"
$ws.Range("E23").Value = "  +0.55%  "

# Row 24
$ws.Range("D24").Value = "9.06"
$ws.Range("E24").Value = "  -3.38%  "

# Row 25
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  -0.14%  "

# Row 26
$ws.Range("D26").Value = "2.50"
$ws.Range("E26").Value = "  -1.51%  "

# Row 27
$ws.Range("D27").Value = "11.30"
$ws.Range("E27").Value = "  -2.94%  "

# Row 28
$ws.Range("E28").Value = "  -1.23%  "

# Row 29
$ws.Range("E29").Value = "  +1.26%  "

# Row 30
$ws.Range("D30").Value = "37.77"
$ws.Range("E30").Value = "  -10.18%  "

# Row 31
$ws.Range("D31").Value = "175.14"
$ws.Range("E31").Value = "  +1.31%  "

# Row 32
$ws.Range("D32").Value = "22.15"
$ws.Range("E32").Value = "  +4.83%  "

# Row 33
$ws.Range("D33").Value = "0.0891"
$ws.Range("E33").Value = "  -3.80%  "

# Row 34
$ws.Range("D34").Value = "5.44"
$ws.Range("E34").Value = "  -4.10%  "

# Row 35
$ws.Range("D35").Value = "0.126"
$ws.Range("E35").Value = "  +0.86%  "

# Row 36
$ws.Range("D36").Value = "0.109"
$ws.Range("E36").Value = "  -5.35%  "

# Row 37
$ws.Range("D37").Value = "4.48"
$ws.Range("E37").Value = "  +1.41%  "

# Row 38
$ws.Range("D38").Value = "0.0357"
$ws.Range("E38").Value = "  -0.07%  "

# Row 39
$ws.Range("D39").Value = "3.27"
$ws.Range("E39").Value = "  -13.20%  "

# Row 40
$ws.Range("D40").Value = "2.36"
$ws.Range("E40").Value = "  +7.59%  "

# Row 41
$ws.Range("D41").Value = "0.237"
$ws.Range("E41").Value = "  -0.55%  "

# Row 42
$ws.Range("D42").Value = "1.37"
$ws.Range("E42").Value = "  +16.72%  "

# Row 43
$ws.Range("D43").Value = "12.11"
$ws.Range("E43").Value = "  -8.26%  "

# Row 44
$ws.Range("D44").Value = "62.20"
$ws.Range("E44").Value = "  -1.85%  "

# Row 45
$ws.Range("D45").Value = "9.03"
$ws.Range("E45").Value = "  +4.99%  "

# Row 46
$ws.Range("D46").Value = "5.28"
$ws.Range("E46").Value = "  -3.91%  "

# Row 47
$ws.Range("E47").Value = "  +1.50%  "

# Row 48
$ws.Range("D48").Value = "100.28"
$ws.Range("E48").Value = "  -2.51%  "

# Row 49
$ws.Range("D49").Value = "1.20"
$ws.Range("E49").Value = "  -0.45%  "

# Row 50
$ws.Range("D50").Value = "2.533.06"
$ws.Range("E50").Value = "  +2.07%  "

# Row 51
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "1.49"
$ws.Range("E51").Value = "  +2.46%  "
